# Apply the edit described by the commit:
#  - update the absPath (now a OneDrive URL) recorded in workbook.xml
#  - insert a new merged/centered "Dates" header row above the existing date row on Sheet1
#  - add a new Sheet2 with Column1/Column2/Column3 headers
#  - leave the selections where the author left them (G8 on Sheet1, D5 on Sheet2)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add Sheet2 (after Sheet1) with the new column headers ---
# (populated first so the shared-string table picks up Column1/2/3 ahead of "Dates")
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Column1"
$ws2.Range("B1").Value = "Column2"
$ws2.Range("C1").Value = "Column3"

# --- Sheet1: push the existing data down one row and add a merged/centered header ---
$ws1.Rows("1:1").Insert()

$ws1.Range("A1").Value = "Dates"
$ws1.Range("A1:D1").HorizontalAlignment = -4108   # xlCenter
$ws1.Range("A1:D1").Merge()

$ws2.Range("D5").Select()
$ws1.Activate()
$ws1.Range("G8").Select()
